# Add the 2018 wealth tax ("Vermögenssteuer", Steuertyp=4) tariff rows for
# the Staatssteuer sheet (ZH, single + married tariffs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    # Jahr, Kanton, Steuertyp, Tariftyp, Einkommen, Steuer, Steuerinkrement, Einkommensinkrement
    @(2018, "ZH", 4, 1,   77000,    0,      0.5, 1000),
    @(2018, "ZH", 4, 1,  308000,  115.5,    1,   1000),
    @(2018, "ZH", 4, 1,  694000,  501.5,    1.5, 1000),
    @(2018, "ZH", 4, 1, 1310000, 1425.5,    2,   1000),
    @(2018, "ZH", 4, 1, 2235000, 3275.5,    2.5, 1000),
    @(2018, "ZH", 4, 1, 3158000, 5583,      3,   1000),
    @(2018, "ZH", 4, 2,  154000,    0,      0.5, 1000),
    @(2018, "ZH", 4, 2,  385000,  115.5,    1,   1000),
    @(2018, "ZH", 4, 2,  770000,  500.5,    1.5, 1000),
    @(2018, "ZH", 4, 2, 1386000, 1424.5,    2,   1000),
    @(2018, "ZH", 4, 2, 2311000, 3274.5,    2.5, 1000),
    @(2018, "ZH", 4, 2, 3235000, 5584.5,    3,   1000)
)

$startRow = 26
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}

# Resize columns to fit the new data (mirrors the column-width shrink seen
# after adding the narrower wealth-tax columns).
$ws.Range("C1").ColumnWidth = 4.666666666666667
$ws.Range("D1").ColumnWidth = 3.8333333333333335
$ws.Range("E1").ColumnWidth = 9
$ws.Range("G1").ColumnWidth = 6.666666666666667
$ws.Range("H1").ColumnWidth = 12.5

# Update the view: scroll down and select the new last numeric cell.
$ws.Range("H28").Select() | Out-Null
